$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the card-detail rows (2-7 and 8-13) into two single cells
# formatted as Python tuple-literal strings, then remove the now-unused
# trailing rows.

$ws.Range("A2").Value = "('Force of Nature', ['{2}{G}{G}{G}{G}', 'Creature — Elemental', 'Trample (This creature can deal excess combat damage to the player or planeswalker it’s attacking.)', 'At the beginning of your upkeep, Force of Nature deals 8 damage to you unless you pay {G}{G}{G}{G}.', '8/8'])"

$ws.Range("A3").Value = "('Шиванский дракон', ['{4}{R}{R}', 'Существо — Дракон', 'Полет (Это существо может быть блокировано только существом со способностью полета).', '{R}: Шиванский дракон получает +1/+0 до конца хода.', '5/5'])"

$ws.Range("A4:A13").ClearContents()
